# Extraction now chooses ideal intrinsic width.
#
# "listing" sheet: column G ("Chosen Intrinsic Width") becomes a constant
#   380 for every data row (2-17) instead of copying column F per row.
# "detail" sheet: column G ("Chosen Intrinsic Width") becomes a constant
#   570 for every data row (2-17); also column A ("Usage") values are
#   rescaled by 1/100 (they were stored as percentages but need to be
#   plain fractions).

$wb = $excel.ActiveWorkbook

$listing = $wb.Worksheets.Item("listing")
$detail  = $wb.Worksheets.Item("detail")

# listing!G2:G17 -> 380 (chosen intrinsic width is now a single fixed value)
$listing.Range("G2:G17").Value2 = 380

# detail!G2:G17 -> 570 (chosen intrinsic width is now a single fixed value)
$detail.Range("G2:G17").Value2 = 570

# detail!A2:A17 -> rescale existing "Usage" fraction by 1/100
for ($r = 2; $r -le 17; $r++) {
    $cell = $detail.Cells.Item($r, 1)
    $old = $cell.Value2
    $cell.Value2 = $old / 100
}
